$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 16 (pushes existing rows 16..56 down to 17..57,
# growing the used range from A1:T56 to A1:T57).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row with this week's Frambuesa record.
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44525
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = "Berries"
$ws.Range("I16").Value = 100101004
$ws.Range("J16").Value = "Frambuesa"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 350
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("Q16").Value = "$/bandeja 2 kilos"
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value = 4000
$ws.Range("T16").Value = 2
